# CaseStudies12Mar2024.pptx — "Add files via upload" edit
#
# Content change: the title on slide 2 ("Title 1" shape) is retitled from
# "Editorial Documentation and Document Management System" to
# "Editorial Pagination and Document Management System" (only the word
# "Documentation" -> "Pagination" changes; surrounding text/whitespace,
# including the trailing tab characters, is preserved as-is).

$p = $ppt.ActivePresentation

$target = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*Editorial Documentation and Document Management System*") {
                $target = $tr
                break
            }
        }
    }
    if ($target -ne $null) { break }
}

if ($target -ne $null) {
    $target.Text = $target.Text.Replace("Editorial Documentation", "Editorial Pagination")
}
